# DOM and Banner author ids added
# - Row 2: fix "Paivi Piirilä" -> "Päivi Piirilä" in the author list, and bump
#   cited_by_count (M2) from 10 to 11.
# - Rows 5-8: the four publication rows got re-ordered (row 5 <-> row 7,
#   row 6 <-> row 8), plus two small data fixes that travel with the rows:
#     * the row that ends up at row 5 (NAMPT/ARDS paper) gets cited_by_count
#       bumped from 0 to 1
#     * the row that ends up at row 8 (Nonequilibrium thermodynamics paper)
#       has its author "Lisa B. Davidson" corrected to "Lisa Davidson"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowRange {
    param($RangeA, $RangeB)

    $rngA = $ws.Range($RangeA)
    $rngB = $ws.Range($RangeB)

    # Capture values before touching anything.
    $valsA = $rngA.Value2
    $valsB = $rngB.Value2

    # Remember the original look (style/format) of each range so the swap
    # below -- which forces Text format to stop Excel from re-typing
    # numeric-looking strings like "0"/"2023" as numbers -- doesn't leave
    # any lasting formatting change behind.
    $styleA = $rngA.Style
    $styleB = $rngB.Style

    $rngA.NumberFormat = "@"
    $rngB.NumberFormat = "@"

    $rngA.Value2 = $valsB
    $rngB.Value2 = $valsA

    $rngA.Style = $styleA
    $rngB.Style = $styleB
}

function Set-TextValue {
    param($CellRef, $NewValue)

    $c = $ws.Range($CellRef)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value2 = $NewValue
    $c.Style = $origStyle
}

# --- Row 2 fixes -----------------------------------------------------
$a2 = $ws.Range("A2").Value2
$a2 = $a2 -replace "Paivi Piirilä", "Päivi Piirilä"
Set-TextValue "A2" $a2
Set-TextValue "M2" "11"

# --- Rows 5-8 re-order -------------------------------------------------
Swap-RowRange "A5:Q5" "A7:Q7"
Swap-RowRange "A6:Q6" "A8:Q8"

# The NAMPT/ARDS paper now lives in row 5 -- bump its citation count.
Set-TextValue "M5" "1"

# The "Nonequilibrium thermodynamics..." paper now lives in row 8 -- fix
# the author name.
$a8 = $ws.Range("A8").Value2
$a8 = $a8 -replace "Lisa B\. Davidson", "Lisa Davidson"
Set-TextValue "A8" $a8
